# Handling partially vanishing saddles
# Updates the "Done / Not Done" status textboxes on slide 1 of the
# FittingEnergy deck.

$p = $ppt.ActivePresentation
$slide = $p.Slides.Item(1)

# ---------------------------------------------------------------
# 1. Remove the two red "(Not Done)" textboxes (ids 11 and 47).
#    Collect them first so deleting doesn't disturb the loop index.
# ---------------------------------------------------------------
$toDelete = @()
for ($i = 1; $i -le $slide.Shapes.Count; $i++) {
    $shp = $slide.Shapes.Item($i)
    if ($shp.Id -eq 11 -or $shp.Id -eq 47) {
        $toDelete += $shp
    }
}
foreach ($shp in $toDelete) {
    $shp.Delete()
}

# ---------------------------------------------------------------
# 2. Move/resize/retext the surviving green "(Done)" textbox
#    (id 46) into "(Done/Kind of)".
# ---------------------------------------------------------------
for ($i = 1; $i -le $slide.Shapes.Count; $i++) {
    $shp = $slide.Shapes.Item($i)
    if ($shp.Id -eq 46) {
        $shp.Left = 808.6966929133858
        $shp.Top = 190.7263779527559
        $shp.Width = 124.62023622047244
        $shp.Height = 29.081259842519685
        $shp.TextFrame.TextRange.Text = "(Done/Kind of)"
        $shp.TextFrame.TextRange.Font.Color.RGB = 0x50B000
    }
}

# ---------------------------------------------------------------
# 3. Add the new green "(Done)" textbox that replaces the old
#    "(Not Done)" box that used to sit at (7886224, 3749052).
# ---------------------------------------------------------------
$newBox1 = $slide.Shapes.AddTextbox(1, 604.1624409448818, 255.45937007874016, 65.15519685039371, 29.081259842519685)
$newBox1.Name = "文本框 1"
$newBox1.TextFrame.WordWrap = 0
$newBox1.TextFrame.AutoSize = 1
$newBox1.Fill.Visible = 0
$newBox1.TextFrame.TextRange.Text = "(Done)"
$newBox1.TextFrame.TextRange.Font.Color.RGB = 0x50B000

# ---------------------------------------------------------------
# 4. Add the new green "(Done)" textbox that replaces the old
#    "(Not Done)" box that used to sit at (7886223, 3234462).
# ---------------------------------------------------------------
$newBox2 = $slide.Shapes.AddTextbox(1, 636.74, 292.87527559055115, 65.15519685039371, 29.081259842519685)
$newBox2.TextFrame.WordWrap = 0
$newBox2.TextFrame.AutoSize = 1
$newBox2.Fill.Visible = 0
$newBox2.TextFrame.TextRange.Text = "(Done)"
$newBox2.TextFrame.TextRange.Font.Color.RGB = 0x50B000
